# Refresh the "Season_Trophies/88" leaderboard sheet with freshly scraped standings.
# Rows 27-114 get new player data (A=Rank, B=ID, C=Name, D=Type/Hall, E=season score),
# and the previous last row (115) is removed, shrinking the sheet from A1:E115 to A1:E114.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data cells in this sheet are stored as text (inline strings), even the purely
# numeric-looking ones (Rank/ID/score). Force the Text number format first so the
# values we assign below are kept as text instead of being coerced into numbers.
$ws.Range("A27:E114").NumberFormat = "@"

# row 27
$ws.Range("A27").Value = "26418"
$ws.Range("B27").Value = "20737010"
$ws.Range("C27").Value = "混着玩..."
$ws.Range("D27").Value = "一馆"
$ws.Range("E27").Value = "4909"

# row 28
$ws.Range("A28").Value = "0"
$ws.Range("B28").Value = "24733875"
$ws.Range("C28").Value = "龍少"
$ws.Range("D28").Value = "一馆"
$ws.Range("E28").Value = "0"

# row 29
$ws.Range("A29").Value = "21649"
$ws.Range("B29").Value = "27484940"
$ws.Range("C29").Value = "66666"
$ws.Range("D29").Value = "一馆"
$ws.Range("E29").Value = "5221"

# row 30
$ws.Range("A30").Value = "16212"
$ws.Range("B30").Value = "31134300"
$ws.Range("C30").Value = "McMaX"
$ws.Range("D30").Value = "一馆"
$ws.Range("E30").Value = "5542"

# row 31
$ws.Range("A31").Value = "16751"
$ws.Range("B31").Value = "31495601"
$ws.Range("C31").Value = "陈晓军"
$ws.Range("D31").Value = "一馆"
$ws.Range("E31").Value = "5505"

# row 32
$ws.Range("A32").Value = "24680"
$ws.Range("B32").Value = "32316256"
$ws.Range("C32").Value = "`"秋の風 ..`""
$ws.Range("D32").Value = "一馆"
$ws.Range("E32").Value = "5015"

# row 33
$ws.Range("A33").Value = "33875"
$ws.Range("B33").Value = "37069173"
$ws.Range("C33").Value = "詹toniii"
$ws.Range("D33").Value = "一馆"
$ws.Range("E33").Value = "4556"

# row 34
$ws.Range("A34").Value = "65450"
$ws.Range("B34").Value = "37861953"
$ws.Range("C34").Value = "`"Durex ๑• . •๑`""
$ws.Range("D34").Value = "一馆"
$ws.Range("E34").Value = "2601"

# row 35
$ws.Range("A35").Value = "19129"
$ws.Range("B35").Value = "38809086"
$ws.Range("C35").Value = "Kouenᶻᵍˣ"
$ws.Range("D35").Value = "一馆"
$ws.Range("E35").Value = "5358"

# row 36
$ws.Range("A36").Value = "31214"
$ws.Range("B36").Value = "38893233"
$ws.Range("C36").Value = "`"快乐 二哈`""
$ws.Range("D36").Value = "一馆"
$ws.Range("E36").Value = "4673"

# row 37
$ws.Range("A37").Value = "21736"
$ws.Range("B37").Value = "38995116"
$ws.Range("C37").Value = "`"Ramesh Pavai Nam`""
$ws.Range("D37").Value = "一馆"
$ws.Range("E37").Value = "5216"

# row 38
$ws.Range("A38").Value = "12903"
$ws.Range("B38").Value = "45967307"
$ws.Range("C38").Value = "Ricky"
$ws.Range("D38").Value = "一馆"
$ws.Range("E38").Value = "5849"

# row 39
$ws.Range("A39").Value = "29387"
$ws.Range("B39").Value = "46289694"
$ws.Range("C39").Value = "㊥Vincent"
$ws.Range("D39").Value = "一馆"
$ws.Range("E39").Value = "4755"

# row 40
$ws.Range("A40").Value = "60627"
$ws.Range("B40").Value = "47146736"
$ws.Range("C40").Value = "`"HK 品瑜`""
$ws.Range("D40").Value = "一馆"
$ws.Range("E40").Value = "2731"

# row 41
$ws.Range("A41").Value = "30612"
$ws.Range("B41").Value = "47459684"
$ws.Range("C41").Value = "㊥阿闹切克闹"
$ws.Range("D41").Value = "一馆"
$ws.Range("E41").Value = "4700"

# row 42
$ws.Range("A42").Value = "46784"
$ws.Range("B42").Value = "48634530"
$ws.Range("C42").Value = "leezhenrui"
$ws.Range("D42").Value = "一馆"
$ws.Range("E42").Value = "3600"

# row 43
$ws.Range("A43").Value = "13419"
$ws.Range("B43").Value = "49043337"
$ws.Range("C43").Value = "FanXiFang1976"
$ws.Range("D43").Value = "一馆"
$ws.Range("E43").Value = "5787"

# row 44
$ws.Range("A44").Value = "11659"
$ws.Range("B44").Value = "49710892"
$ws.Range("C44").Value = "MMMMMMM"
$ws.Range("D44").Value = "一馆"
$ws.Range("E44").Value = "6019"

# row 45
$ws.Range("A45").Value = "38893"
$ws.Range("B45").Value = "50837459"
$ws.Range("C45").Value = "NINE日"
$ws.Range("D45").Value = "一馆"
$ws.Range("E45").Value = "4271"

# row 46
$ws.Range("A46").Value = "42440"
$ws.Range("B46").Value = "52997727"
$ws.Range("C46").Value = "larios"
$ws.Range("D46").Value = "一馆"
$ws.Range("E46").Value = "4050"

# row 47
$ws.Range("A47").Value = "10747"
$ws.Range("B47").Value = "53060417"
$ws.Range("C47").Value = "㊥老纳信耶稣"
$ws.Range("D47").Value = "一馆"
$ws.Range("E47").Value = "6077"

# row 48
$ws.Range("A48").Value = "13877"
$ws.Range("B48").Value = "53520939"
$ws.Range("C48").Value = "㊥虎哥tiger"
$ws.Range("D48").Value = "一馆"
$ws.Range("E48").Value = "5739"

# row 49
$ws.Range("A49").Value = "21786"
$ws.Range("B49").Value = "54085771"
$ws.Range("C49").Value = "㊥Matthieu"
$ws.Range("D49").Value = "一馆"
$ws.Range("E49").Value = "5212"

# row 50
$ws.Range("A50").Value = "16162"
$ws.Range("B50").Value = "54698813"
$ws.Range("C50").Value = "閃亮唐老鴨"
$ws.Range("D50").Value = "一馆"
$ws.Range("E50").Value = "5545"

# row 51
$ws.Range("A51").Value = "24104"
$ws.Range("B51").Value = "54778421"
$ws.Range("C51").Value = "Emma"
$ws.Range("D51").Value = "一馆"
$ws.Range("E51").Value = "5049"

# row 52
$ws.Range("A52").Value = "16439"
$ws.Range("B52").Value = "55317038"
$ws.Range("C52").Value = "necman12345"
$ws.Range("D52").Value = "一馆"
$ws.Range("E52").Value = "5526"

# row 53
$ws.Range("A53").Value = "41082"
$ws.Range("B53").Value = "55634661"
$ws.Range("C53").Value = "Opalus"
$ws.Range("D53").Value = "一馆"
$ws.Range("E53").Value = "4145"

# row 54
$ws.Range("A54").Value = "21126"
$ws.Range("B54").Value = "55769051"
$ws.Range("C54").Value = "㊥叮叮当."
$ws.Range("D54").Value = "一馆"
$ws.Range("E54").Value = "5247"

# row 55
$ws.Range("A55").Value = "25741"
$ws.Range("B55").Value = "55860890"
$ws.Range("C55").Value = "㊥Ethan"
$ws.Range("D55").Value = "一馆"
$ws.Range("E55").Value = "4947"

# row 56
$ws.Range("A56").Value = "15732"
$ws.Range("B56").Value = "56133764"
$ws.Range("C56").Value = "ustcarter"
$ws.Range("D56").Value = "一馆"
$ws.Range("E56").Value = "5577"

# row 57
$ws.Range("A57").Value = "34329"
$ws.Range("B57").Value = "56379103"
$ws.Range("C57").Value = "Globalking"
$ws.Range("D57").Value = "一馆"
$ws.Range("E57").Value = "4535"

# row 58
$ws.Range("A58").Value = "35187"
$ws.Range("B58").Value = "56573048"
$ws.Range("C58").Value = "Xiaotian"
$ws.Range("D58").Value = "一馆"
$ws.Range("E58").Value = "4492"

# row 59
$ws.Range("A59").Value = "35943"
$ws.Range("B59").Value = "56585361"
$ws.Range("C59").Value = "`"㊥ go策划我要ali`""
$ws.Range("D59").Value = "一馆"
$ws.Range("E59").Value = "4449"

# row 60
$ws.Range("A60").Value = "32956"
$ws.Range("B60").Value = "56732705"
$ws.Range("C60").Value = "时间温柔皆遗憾"
$ws.Range("D60").Value = "一馆"
$ws.Range("E60").Value = "4594"

# row 61
$ws.Range("A61").Value = "41857"
$ws.Range("B61").Value = "57813281"
$ws.Range("C61").Value = "XAUEN"
$ws.Range("D61").Value = "一馆"
$ws.Range("E61").Value = "4095"

# row 62
$ws.Range("A62").Value = "40031"
$ws.Range("B62").Value = "58203298"
$ws.Range("C62").Value = "权旨qua"
$ws.Range("D62").Value = "一馆"
$ws.Range("E62").Value = "4204"

# row 63
$ws.Range("A63").Value = "39083"
$ws.Range("B63").Value = "58408326"
$ws.Range("C63").Value = "`"Killer Bee`""
$ws.Range("D63").Value = "一馆"
$ws.Range("E63").Value = "4259"

# row 64
$ws.Range("A64").Value = "32952"
$ws.Range("B64").Value = "58839983"
$ws.Range("C64").Value = "每逢佳节胖六斤"
$ws.Range("D64").Value = "一馆"
$ws.Range("E64").Value = "4594"

# row 65
$ws.Range("A65").Value = "43036"
$ws.Range("B65").Value = "59020292"
$ws.Range("C65").Value = "Sharnoth"
$ws.Range("D65").Value = "一馆"
$ws.Range("E65").Value = "4007"

# row 66
$ws.Range("A66").Value = "0"
$ws.Range("B66").Value = "1222440"
$ws.Range("C66").Value = "`"Sneaky Ninja Panda`""
$ws.Range("D66").Value = "二馆"
$ws.Range("E66").Value = "0"

# row 67
$ws.Range("A67").Value = "54204"
$ws.Range("B67").Value = "3391765"
$ws.Range("C67").Value = "马er"
$ws.Range("D67").Value = "二馆"
$ws.Range("E67").Value = "2993"

# row 68
$ws.Range("A68").Value = "71285"
$ws.Range("B68").Value = "9718882"
$ws.Range("C68").Value = "小霸王2021"
$ws.Range("D68").Value = "二馆"
$ws.Range("E68").Value = "2502"

# row 69
$ws.Range("A69").Value = "36132"
$ws.Range("B69").Value = "11645391"
$ws.Range("C69").Value = "`"omar omar`""
$ws.Range("D69").Value = "二馆"
$ws.Range("E69").Value = "4437"

# row 70
$ws.Range("A70").Value = "92306"
$ws.Range("B70").Value = "15436348"
$ws.Range("C70").Value = "Lucas"
$ws.Range("D70").Value = "二馆"
$ws.Range("E70").Value = "1972"

# row 71
$ws.Range("A71").Value = "0"
$ws.Range("B71").Value = "20372140"
$ws.Range("C71").Value = "人山即是仙"
$ws.Range("D71").Value = "二馆"
$ws.Range("E71").Value = "0"

# row 72
$ws.Range("A72").Value = "0"
$ws.Range("B72").Value = "38994054"
$ws.Range("C72").Value = "chengnan"
$ws.Range("D72").Value = "二馆"
$ws.Range("E72").Value = "0"

# row 73
$ws.Range("A73").Value = "54695"
$ws.Range("B73").Value = "41837764"
$ws.Range("C73").Value = "好风光会长"
$ws.Range("D73").Value = "二馆"
$ws.Range("E73").Value = "2964"

# row 74
$ws.Range("A74").Value = "59977"
$ws.Range("B74").Value = "43281368"
$ws.Range("C74").Value = "xhs2763"
$ws.Range("D74").Value = "二馆"
$ws.Range("E74").Value = "2752"

# row 75
$ws.Range("A75").Value = "0"
$ws.Range("B75").Value = "44378757"
$ws.Range("C75").Value = "`"NᵉᵗʰᵉʳDʳⁱᶠᵗᵉʳ ㊥`""
$ws.Range("D75").Value = "二馆"
$ws.Range("E75").Value = "0"

# row 76
$ws.Range("A76").Value = "31659"
$ws.Range("B76").Value = "44708798"
$ws.Range("C76").Value = "`"㊥ mythgod`""
$ws.Range("D76").Value = "二馆"
$ws.Range("E76").Value = "4652"

# row 77
$ws.Range("A77").Value = "46612"
$ws.Range("B77").Value = "47430231"
$ws.Range("C77").Value = "Kentantrino"
$ws.Range("D77").Value = "二馆"
$ws.Range("E77").Value = "3621"

# row 78
$ws.Range("A78").Value = "0"
$ws.Range("B78").Value = "48738257"
$ws.Range("C78").Value = "死亡洲际跳蛋"
$ws.Range("D78").Value = "二馆"
$ws.Range("E78").Value = "0"

# row 79
$ws.Range("A79").Value = "0"
$ws.Range("B79").Value = "49000199"
$ws.Range("C79").Value = "SlipperyForester5672"
$ws.Range("D79").Value = "二馆"
$ws.Range("E79").Value = "1225"

# row 80
$ws.Range("A80").Value = "0"
$ws.Range("B80").Value = "54941706"
$ws.Range("C80").Value = "AlexMenjivar20"
$ws.Range("D80").Value = "二馆"
$ws.Range("E80").Value = "1496"

# row 81
$ws.Range("A81").Value = "51979"
$ws.Range("B81").Value = "55499394"
$ws.Range("C81").Value = "Player-55499394"
$ws.Range("D81").Value = "二馆"
$ws.Range("E81").Value = "3126"

# row 82
$ws.Range("A82").Value = "0"
$ws.Range("B82").Value = "55810157"
$ws.Range("C82").Value = "Beard"
$ws.Range("D82").Value = "二馆"
$ws.Range("E82").Value = "0"

# row 83
$ws.Range("A83").Value = "0"
$ws.Range("B83").Value = "56700848"
$ws.Range("C83").Value = "工口漫画老师"
$ws.Range("D83").Value = "二馆"
$ws.Range("E83").Value = "0"

# row 84
$ws.Range("A84").Value = "0"
$ws.Range("B84").Value = "57219176"
$ws.Range("C84").Value = "青莲道人"
$ws.Range("D84").Value = "二馆"
$ws.Range("E84").Value = "1525"

# row 85
$ws.Range("A85").Value = "0"
$ws.Range("B85").Value = "57556179"
$ws.Range("C85").Value = "特战新生代英雄"
$ws.Range("D85").Value = "二馆"
$ws.Range("E85").Value = "0"

# row 86
$ws.Range("A86").Value = "0"
$ws.Range("B86").Value = "58340439"
$ws.Range("C86").Value = "70qilin"
$ws.Range("D86").Value = "二馆"
$ws.Range("E86").Value = "0"

# row 87
$ws.Range("A87").Value = "0"
$ws.Range("B87").Value = "58615925"
$ws.Range("C87").Value = "齐天的大圣"
$ws.Range("D87").Value = "二馆"
$ws.Range("E87").Value = "0"

# row 88
$ws.Range("A88").Value = "0"
$ws.Range("B88").Value = "58641574"
$ws.Range("C88").Value = "Player-58641574鱼"
$ws.Range("D88").Value = "二馆"
$ws.Range("E88").Value = "0"

# row 89
$ws.Range("A89").Value = "0"
$ws.Range("B89").Value = "58743790"
$ws.Range("C89").Value = "Ma"
$ws.Range("D89").Value = "二馆"
$ws.Range("E89").Value = "0"

# row 90
$ws.Range("A90").Value = "0"
$ws.Range("B90").Value = "15695258"
$ws.Range("C90").Value = "Player-15695258"
$ws.Range("D90").Value = "三馆"
$ws.Range("E90").Value = "1000"

# row 91
$ws.Range("A91").Value = "0"
$ws.Range("B91").Value = "29355299"
$ws.Range("C91").Value = "Player-29355299"
$ws.Range("D91").Value = "三馆"
$ws.Range("E91").Value = "1000"

# row 92
$ws.Range("A92").Value = "0"
$ws.Range("B92").Value = "41231396"
$ws.Range("C92").Value = "ollsthebro"
$ws.Range("D92").Value = "三馆"
$ws.Range("E92").Value = "0"

# row 93
$ws.Range("A93").Value = "0"
$ws.Range("B93").Value = "47622456"
$ws.Range("C93").Value = "伊恩"
$ws.Range("D93").Value = "三馆"
$ws.Range("E93").Value = "0"

# row 94
$ws.Range("A94").Value = "0"
$ws.Range("B94").Value = "49553719"
$ws.Range("C94").Value = "`"Oreo Captain Sir`""
$ws.Range("D94").Value = "三馆"
$ws.Range("E94").Value = "0"

# row 95
$ws.Range("A95").Value = "0"
$ws.Range("B95").Value = "55745105"
$ws.Range("C95").Value = "eldeniz"
$ws.Range("D95").Value = "三馆"
$ws.Range("E95").Value = "1000"

# row 96
$ws.Range("A96").Value = "45837"
$ws.Range("B96").Value = "56241637"
$ws.Range("C96").Value = "Player-14day"
$ws.Range("D96").Value = "三馆"
$ws.Range("E96").Value = "3713"

# row 97
$ws.Range("A97").Value = "0"
$ws.Range("B97").Value = "58174442"
$ws.Range("C97").Value = "Player-58174442"
$ws.Range("D97").Value = "三馆"
$ws.Range("E97").Value = "1000"

# row 98
$ws.Range("A98").Value = "0"
$ws.Range("B98").Value = "58572199"
$ws.Range("C98").Value = "你干嘛～哎呦～"
$ws.Range("D98").Value = "三馆"
$ws.Range("E98").Value = "0"

# row 99
$ws.Range("A99").Value = "0"
$ws.Range("B99").Value = "58671339"
$ws.Range("C99").Value = "`"quang pro`""
$ws.Range("D99").Value = "三馆"
$ws.Range("E99").Value = "0"

# row 100
$ws.Range("A100").Value = "0"
$ws.Range("B100").Value = "58766144"
$ws.Range("C100").Value = "EquablePrecedence38"
$ws.Range("D100").Value = "三馆"
$ws.Range("E100").Value = "0"

# row 101
$ws.Range("A101").Value = "0"
$ws.Range("B101").Value = "58910668"
$ws.Range("C101").Value = "BrittleAuthor33"
$ws.Range("D101").Value = "三馆"
$ws.Range("E101").Value = "0"

# row 102
$ws.Range("A102").Value = "0"
$ws.Range("B102").Value = "59081265"
$ws.Range("C102").Value = "爬楼梯"
$ws.Range("D102").Value = "三馆"
$ws.Range("E102").Value = "0"

# row 103
$ws.Range("A103").Value = "0"
$ws.Range("B103").Value = "59082827"
$ws.Range("C103").Value = "Player-59082827"
$ws.Range("D103").Value = "三馆"
$ws.Range("E103").Value = "0"

# row 104
$ws.Range("A104").Value = "0"
$ws.Range("B104").Value = "59106471"
$ws.Range("C104").Value = "anime"
$ws.Range("D104").Value = "三馆"
$ws.Range("E104").Value = "1498"

# row 105
$ws.Range("A105").Value = "0"
$ws.Range("B105").Value = "59112086"
$ws.Range("C105").Value = "sigma"
$ws.Range("D105").Value = "三馆"
$ws.Range("E105").Value = "0"

# row 106
$ws.Range("A106").Value = "71461"
$ws.Range("B106").Value = "6010122"
$ws.Range("C106").Value = "`"Edward Peng`""
$ws.Range("D106").Value = "Chinese"
$ws.Range("E106").Value = "0"

# row 107
$ws.Range("A107").Value = "0"
$ws.Range("B107").Value = "8850180"
$ws.Range("C107").Value = "30624300"
$ws.Range("D107").Value = "Chinese"
$ws.Range("E107").Value = "0"

# row 108
$ws.Range("A108").Value = "0"
$ws.Range("B108").Value = "9195340"
$ws.Range("C108").Value = "Namllllllik"
$ws.Range("D108").Value = "Chinese"
$ws.Range("E108").Value = "0"

# row 109
$ws.Range("A109").Value = "69265"
$ws.Range("B109").Value = "9913517"
$ws.Range("C109").Value = "`"Kenny Chan`""
$ws.Range("D109").Value = "Chinese"
$ws.Range("E109").Value = "0"

# row 110
$ws.Range("A110").Value = "0"
$ws.Range("B110").Value = "10636651"
$ws.Range("C110").Value = "`"Ismail Aflou`""
$ws.Range("D110").Value = "Chinese"
$ws.Range("E110").Value = "0"

# row 111
$ws.Range("A111").Value = "0"
$ws.Range("B111").Value = "12648101"
$ws.Range("C111").Value = "`"player 198827`""
$ws.Range("D111").Value = "Chinese"
$ws.Range("E111").Value = "0"

# row 112
$ws.Range("A112").Value = "71670"
$ws.Range("B112").Value = "15755724"
$ws.Range("C112").Value = "`"Last Good`""
$ws.Range("D112").Value = "Chinese"
$ws.Range("E112").Value = "0"

# row 113
$ws.Range("A113").Value = "0"
$ws.Range("B113").Value = "28624723"
$ws.Range("C113").Value = "`"Woody Shade`""
$ws.Range("D113").Value = "Chinese"
$ws.Range("E113").Value = "0"

# row 114
$ws.Range("A114").Value = "64399"
$ws.Range("B114").Value = "41848598"
$ws.Range("C114").Value = "国家一级保护沙雕"
$ws.Range("D114").Value = "Chinese"
$ws.Range("E114").Value = "0"

# The old row 115 is gone in the new data; delete it so the sheet ends at row 114
# (this also updates the sheet dimension to A1:E114 automatically).
$ws.Rows(115).Delete()
